$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 29.223446
$ws.Cells.Item(2, 8).Value = 87.670338
$ws.Cells.Item(2, 9).Value = 0.0169041244192178
$ws.Cells.Item(2, 10).Value = 0.0169041244192178
$ws.Cells.Item(2, 13).Value = 2.113523666666667
$ws.Cells.Item(2, 14).Value = 6.340571000000001
$ws.Cells.Item(2, 15).Value = 0.2651220308693004
$ws.Cells.Item(2, 16).Value = 0.2651220308693004
$ws.Cells.Item(2, 17).Value = 61.76444474255534
$ws.Cells.Item(2, 18).Value = 555.8800026829981
$ws.Cells.Item(2, 19).Value = 0.004481655796090356
$ws.Cells.Item(2, 20).Value = 0.004481655796090356

$ws.Cells.Item(3, 7).Value = 29.223446
$ws.Cells.Item(3, 8).Value = 87.670338
$ws.Cells.Item(3, 9).Value = 0.0169041244192178
$ws.Cells.Item(3, 10).Value = 0.0169041244192178
$ws.Cells.Item(3, 15).Value = 0.2869289465860668
$ws.Cells.Item(3, 16).Value = 0.2869289465860668
$ws.Cells.Item(3, 17).Value = 66.844716783236
$ws.Cells.Item(3, 18).Value = 601.6024510491239
$ws.Cells.Item(3, 19).Value = 0.004850282612565972
$ws.Cells.Item(3, 20).Value = 0.004850282612565972

$ws.Cells.Item(4, 7).Value = 29.223446
$ws.Cells.Item(4, 8).Value = 87.670338
$ws.Cells.Item(4, 9).Value = 0.0169041244192178
$ws.Cells.Item(4, 10).Value = 0.0169041244192178
$ws.Cells.Item(4, 13).Value = 1.164746666666667
$ws.Cells.Item(4, 14).Value = 3.49424
$ws.Cells.Item(4, 15).Value = 0.1461067158059967
$ws.Cells.Item(4, 16).Value = 0.1461067158059966
$ws.Cells.Item(4, 17).Value = 34.03791131701333
$ws.Cells.Item(4, 18).Value = 306.34120185312
$ws.Cells.Item(4, 19).Value = 0.002469806102467864
$ws.Cells.Item(4, 20).Value = 0.002469806102467863

$ws.Cells.Item(5, 7).Value = 29.223446
$ws.Cells.Item(5, 8).Value = 87.670338
$ws.Cells.Item(5, 9).Value = 0.0169041244192178
$ws.Cells.Item(5, 10).Value = 0.0169041244192178
$ws.Cells.Item(5, 13).Value = 2.406253666666667
$ws.Cells.Item(5, 14).Value = 7.218761
$ws.Cells.Item(5, 15).Value = 0.3018423067386362
$ws.Cells.Item(5, 16).Value = 0.3018423067386362
$ws.Cells.Item(5, 17).Value = 70.31902409013533
$ws.Cells.Item(5, 18).Value = 632.871216811218
$ws.Cells.Item(5, 19).Value = 0.00510237990809361
$ws.Cells.Item(5, 20).Value = 0.00510237990809361

$ws.Cells.Item(6, 9).Value = 0.9471112884046843
$ws.Cells.Item(6, 10).Value = 0.9471112884046842
$ws.Cells.Item(6, 13).Value = 2.113523666666667
$ws.Cells.Item(6, 14).Value = 6.340571000000001
$ws.Cells.Item(6, 15).Value = 0.2651220308693004
$ws.Cells.Item(6, 16).Value = 0.2651220308693004
$ws.Cells.Item(6, 17).Value = 3460.563906594126
$ws.Cells.Item(6, 18).Value = 31145.07515934713
$ws.Cells.Item(6, 19).Value = 0.2511000682410895
$ws.Cells.Item(6, 20).Value = 0.2511000682410895

$ws.Cells.Item(7, 9).Value = 0.9471112884046843
$ws.Cells.Item(7, 10).Value = 0.9471112884046842
$ws.Cells.Item(7, 15).Value = 0.2869289465860668
$ws.Cells.Item(7, 16).Value = 0.2869289465860668
$ws.Cells.Item(7, 19).Value = 0.2717536442817285
$ws.Cells.Item(7, 20).Value = 0.2717536442817285

$ws.Cells.Item(8, 9).Value = 0.9471112884046843
$ws.Cells.Item(8, 10).Value = 0.9471112884046842
$ws.Cells.Item(8, 13).Value = 1.164746666666667
$ws.Cells.Item(8, 14).Value = 3.49424
$ws.Cells.Item(8, 15).Value = 0.1461067158059967
$ws.Cells.Item(8, 16).Value = 0.1461067158059966
$ws.Cells.Item(8, 17).Value = 1907.090201336355
$ws.Cells.Item(8, 18).Value = 17163.8118120272
$ws.Cells.Item(8, 19).Value = 0.1383793198515945
$ws.Cells.Item(8, 20).Value = 0.1383793198515945

$ws.Cells.Item(9, 9).Value = 0.9471112884046843
$ws.Cells.Item(9, 10).Value = 0.9471112884046842
$ws.Cells.Item(9, 13).Value = 2.406253666666667
$ws.Cells.Item(9, 14).Value = 7.218761
$ws.Cells.Item(9, 15).Value = 0.3018423067386362
$ws.Cells.Item(9, 16).Value = 0.3018423067386362
$ws.Cells.Item(9, 17).Value = 3939.863423488092
$ws.Cells.Item(9, 18).Value = 35458.77081139283
$ws.Cells.Item(9, 19).Value = 0.2858782560302717
$ws.Cells.Item(9, 20).Value = 0.2858782560302716

$ws.Cells.Item(10, 7).Value = 37.39212666666667
$ws.Cells.Item(10, 8).Value = 112.17638
$ws.Cells.Item(10, 9).Value = 0.02162924801792661
$ws.Cells.Item(10, 10).Value = 0.0216292480179266
$ws.Cells.Item(10, 13).Value = 2.113523666666667
$ws.Cells.Item(10, 14).Value = 6.340571000000001
$ws.Cells.Item(10, 15).Value = 0.2651220308693004
$ws.Cells.Item(10, 16).Value = 0.2651220308693004
$ws.Cells.Item(10, 17).Value = 79.0291446569978
$ws.Cells.Item(10, 18).Value = 711.2623019129801
$ws.Cells.Item(10, 19).Value = 0.005734390160688492
$ws.Cells.Item(10, 20).Value = 0.005734390160688491

$ws.Cells.Item(11, 7).Value = 37.39212666666667
$ws.Cells.Item(11, 8).Value = 112.17638
$ws.Cells.Item(11, 9).Value = 0.02162924801792661
$ws.Cells.Item(11, 10).Value = 0.0216292480179266
$ws.Cells.Item(11, 15).Value = 0.2869289465860668
$ws.Cells.Item(11, 16).Value = 0.2869289465860668
$ws.Cells.Item(11, 17).Value = 85.52947920502668
$ws.Cells.Item(11, 18).Value = 769.76531284524
$ws.Cells.Item(11, 19).Value = 0.006206057349232455
$ws.Cells.Item(11, 20).Value = 0.006206057349232453

$ws.Cells.Item(12, 7).Value = 37.39212666666667
$ws.Cells.Item(12, 8).Value = 112.17638
$ws.Cells.Item(12, 9).Value = 0.02162924801792661
$ws.Cells.Item(12, 10).Value = 0.0216292480179266
$ws.Cells.Item(12, 13).Value = 1.164746666666667
$ws.Cells.Item(12, 14).Value = 3.49424
$ws.Cells.Item(12, 15).Value = 0.1461067158059967
$ws.Cells.Item(12, 16).Value = 0.1461067158059966
$ws.Cells.Item(12, 17).Value = 43.55235489457778
$ws.Cells.Item(12, 18).Value = 391.9711940512
$ws.Cells.Item(12, 19).Value = 0.003160178393252619
$ws.Cells.Item(12, 20).Value = 0.003160178393252618

$ws.Cells.Item(13, 7).Value = 37.39212666666667
$ws.Cells.Item(13, 8).Value = 112.17638
$ws.Cells.Item(13, 9).Value = 0.02162924801792661
$ws.Cells.Item(13, 10).Value = 0.0216292480179266
$ws.Cells.Item(13, 13).Value = 2.406253666666667
$ws.Cells.Item(13, 14).Value = 7.218761
$ws.Cells.Item(13, 15).Value = 0.3018423067386362
$ws.Cells.Item(13, 16).Value = 0.3018423067386362
$ws.Cells.Item(13, 17).Value = 89.97494189613113
$ws.Cells.Item(13, 18).Value = 809.77447706518
$ws.Cells.Item(13, 19).Value = 0.006528622114753043
$ws.Cells.Item(13, 20).Value = 0.00652862211475304

$ws.Cells.Item(14, 7).Value = 24.817167
$ws.Cells.Item(14, 8).Value = 74.45150100000001
$ws.Cells.Item(14, 9).Value = 0.01435533915817136
$ws.Cells.Item(14, 10).Value = 0.01435533915817136
$ws.Cells.Item(14, 13).Value = 2.113523666666667
$ws.Cells.Item(14, 14).Value = 6.340571000000001
$ws.Cells.Item(14, 15).Value = 0.2651220308693004
$ws.Cells.Item(14, 16).Value = 0.2651220308693004
$ws.Cells.Item(14, 17).Value = 52.45166979411901
$ws.Cells.Item(14, 18).Value = 472.0650281470711
$ws.Cells.Item(14, 19).Value = 0.003805916671431984
$ws.Cells.Item(14, 20).Value = 0.003805916671431984

$ws.Cells.Item(15, 7).Value = 24.817167
$ws.Cells.Item(15, 8).Value = 74.45150100000001
$ws.Cells.Item(15, 9).Value = 0.01435533915817136
$ws.Cells.Item(15, 10).Value = 0.01435533915817136
$ws.Cells.Item(15, 15).Value = 0.2869289465860668
$ws.Cells.Item(15, 16).Value = 0.2869289465860668
$ws.Cells.Item(15, 17).Value = 56.76594401212201
$ws.Cells.Item(15, 18).Value = 510.893496109098
$ws.Cells.Item(15, 19).Value = 0.004118962342539823
$ws.Cells.Item(15, 20).Value = 0.004118962342539822

$ws.Cells.Item(16, 7).Value = 24.817167
$ws.Cells.Item(16, 8).Value = 74.45150100000001
$ws.Cells.Item(16, 9).Value = 0.01435533915817136
$ws.Cells.Item(16, 10).Value = 0.01435533915817136
$ws.Cells.Item(16, 13).Value = 1.164746666666667
$ws.Cells.Item(16, 14).Value = 3.49424
$ws.Cells.Item(16, 15).Value = 0.1461067158059967
$ws.Cells.Item(16, 16).Value = 0.1461067158059966
$ws.Cells.Item(16, 17).Value = 28.90571253936
$ws.Cells.Item(16, 18).Value = 260.15141285424
$ws.Cells.Item(16, 19).Value = 0.002097411458681638
$ws.Cells.Item(16, 20).Value = 0.002097411458681638

$ws.Cells.Item(17, 7).Value = 24.817167
$ws.Cells.Item(17, 8).Value = 74.45150100000001
$ws.Cells.Item(17, 9).Value = 0.01435533915817136
$ws.Cells.Item(17, 10).Value = 0.01435533915817136
$ws.Cells.Item(17, 13).Value = 2.406253666666667
$ws.Cells.Item(17, 14).Value = 7.218761
$ws.Cells.Item(17, 15).Value = 0.3018423067386362
$ws.Cells.Item(17, 16).Value = 0.3018423067386362
$ws.Cells.Item(17, 17).Value = 59.71639909002901
$ws.Cells.Item(17, 18).Value = 537.4475918102611
$ws.Cells.Item(17, 19).Value = 0.004333048685517916
$ws.Cells.Item(17, 20).Value = 0.004333048685517914
